$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A entries in the new rows to be stored as text (not auto-converted to dates)
$dateRange = $ws.Range("A41:A52")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(41, 1).Value = "2025-02-14"
$ws.Cells.Item(41, 2).Value = "sleep"
$ws.Cells.Item(41, 3).Value = $true
$ws.Cells.Item(41, 4).Value = $false

$ws.Cells.Item(42, 1).Value = "2025-02-14"
$ws.Cells.Item(42, 2).Value = "activity"
$ws.Cells.Item(42, 3).Value = $true
$ws.Cells.Item(42, 4).Value = $true

$ws.Cells.Item(43, 1).Value = "2025-02-14"
$ws.Cells.Item(43, 2).Value = "weekly_activity"
$ws.Cells.Item(43, 3).Value = $true
$ws.Cells.Item(43, 4).Value = $false

$ws.Cells.Item(44, 1).Value = "2025-02-15"
$ws.Cells.Item(44, 2).Value = "sleep"
$ws.Cells.Item(44, 3).Value = $true
$ws.Cells.Item(44, 4).Value = $true

$ws.Cells.Item(45, 1).Value = "2025-02-15"
$ws.Cells.Item(45, 2).Value = "activity"
$ws.Cells.Item(45, 3).Value = $true
$ws.Cells.Item(45, 4).Value = $true

$ws.Cells.Item(46, 1).Value = "2025-02-15"
$ws.Cells.Item(46, 2).Value = "weekly_activity"
$ws.Cells.Item(46, 3).Value = $false
$ws.Cells.Item(46, 4).Value = $false

$ws.Cells.Item(47, 1).Value = "2025-02-16"
$ws.Cells.Item(47, 2).Value = "sleep"
$ws.Cells.Item(47, 3).Value = $true
$ws.Cells.Item(47, 4).Value = $true

$ws.Cells.Item(48, 1).Value = "2025-02-16"
$ws.Cells.Item(48, 2).Value = "activity"
$ws.Cells.Item(48, 3).Value = $true
$ws.Cells.Item(48, 4).Value = $true

$ws.Cells.Item(49, 1).Value = "2025-02-16"
$ws.Cells.Item(49, 2).Value = "weekly_activity"
$ws.Cells.Item(49, 3).Value = $false
$ws.Cells.Item(49, 4).Value = $false

$ws.Cells.Item(50, 1).Value = "2025-02-17"
$ws.Cells.Item(50, 2).Value = "sleep"
$ws.Cells.Item(50, 3).Value = $false
$ws.Cells.Item(50, 4).Value = $false

$ws.Cells.Item(51, 1).Value = "2025-02-17"
$ws.Cells.Item(51, 2).Value = "activity"
$ws.Cells.Item(51, 3).Value = $true
$ws.Cells.Item(51, 4).Value = $true

$ws.Cells.Item(52, 1).Value = "2025-02-17"
$ws.Cells.Item(52, 2).Value = "weekly_activity"
$ws.Cells.Item(52, 3).Value = $false
$ws.Cells.Item(52, 4).Value = $false

# Reset style on the date column back to Normal/default so no explicit style index is retained
$dateRange.Style = "Normal"
